$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (e.g. "67.683.13", "0.998")
# that must stay plain text, matching the original inlineStr cells. Force
# the cell to Text format before writing, then reset the style back to
# Normal so we do not leave a stray number-format style behind.
function Set-TextValue($rangeAddr, $val) {
    $c = $ws.Range($rangeAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "67.683.13"
$ws.Range("E2").Value = "  -2.41%  "
Set-TextValue "D3" "3.490.98"
$ws.Range("E3").Value = "  -4.02%  "
Set-TextValue "D4" "0.998"
$ws.Range("E4").Value = "  -0.04%  "
Set-TextValue "D5" "605.79"
$ws.Range("E5").Value = "  -3.27%  "
Set-TextValue "D6" "149.96"
$ws.Range("E6").Value = "  -5.54%  "
Set-TextValue "D7" "3.487.11"
$ws.Range("E7").Value = "  -3.96%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -1.43%  "
Set-TextValue "D10" "0.143"
$ws.Range("E10").Value = "  -0.89%  "
Set-TextValue "D11" "7.52"
$ws.Range("E11").Value = "  +4.55%  "
Set-TextValue "D12" "0.429"
$ws.Range("E12").Value = "  -2.28%  "
Set-TextValue "D13" "0.0000214"
$ws.Range("E13").Value = "  -4.10%  "
Set-TextValue "D14" "31.87"
$ws.Range("E14").Value = "  -4.10%  "
Set-TextValue "D15" "4.076.09"
$ws.Range("E15").Value = "  -4.02%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D16" "3.485.78"
$ws.Range("E16").Value = "  -3.86%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D17" "67.585.80"
$ws.Range("E17").Value = "  -2.71%  "
Set-TextValue "D18" "0.116"
$ws.Range("E18").Value = "  -0.27%  "
Set-TextValue "D19" "6.48"
$ws.Range("E19").Value = "  -1.75%  "
Set-TextValue "D20" "15.31"
$ws.Range("E20").Value = "  -3.30%  "
Set-TextValue "D21" "9.95"
$ws.Range("E21").Value = "  -2.51%  "
Set-TextValue "D22" "443.89"
$ws.Range("E22").Value = "  -3.68%  "
Set-TextValue "D23" "0.623"
$ws.Range("E23").Value = "  -2.89%  "
Set-TextValue "D24" "78.83"
$ws.Range("E24").Value = "  +0.71%  "
Set-TextValue "D25" "3.630.83"
$ws.Range("E25").Value = "  -3.84%  "
$ws.Range("E26").Value = "  -0.31%  "
Set-TextValue "D27" "0.0000123"
$ws.Range("E27").Value = "  -9.53%  "
Set-TextValue "D28" "8.66"
$ws.Range("E28").Value = "  -5.69%  "
$ws.Range("E29").Value = "  -5.74%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D30" "2.50"
$ws.Range("E30").Value = "  -4.49%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D31" "1.65"
$ws.Range("E31").Value = "  -4.90%  "
Set-TextValue "D32" "0.168"
$ws.Range("E32").Value = "  -3.52%  "
$ws.Range("E33").Value = "  +0.10%  "
Set-TextValue "D34" "25.53"
$ws.Range("E34").Value = "  -3.61%  "
Set-TextValue "D35" "6.16"
$ws.Range("E35").Value = "  -6.56%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextValue "D36" "3.480.55"
$ws.Range("E36").Value = "  -3.87%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D37" "1.84"
$ws.Range("E37").Value = "  -5.76%  "
Set-TextValue "D38" "7.97"
$ws.Range("E38").Value = "  -4.59%  "
$ws.Range("E39").Value = "  -0.03%  "
Set-TextValue "D40" "2.29"
$ws.Range("E40").Value = "  -3.13%  "
$ws.Range("E41").Value = "  +0.31%  "
Set-TextValue "D42" "174.90"
$ws.Range("E42").Value = "  +0.38%  "
Set-TextValue "D43" "0.0899"
$ws.Range("E43").Value = "  -2.92%  "
Set-TextValue "D44" "5.40"
$ws.Range("E44").Value = "  -4.50%  "
Set-TextValue "D45" "0.896"
$ws.Range("E45").Value = "  -2.02%  "
Set-TextValue "D46" "30.16"
$ws.Range("E46").Value = "  -2.48%  "
Set-TextValue "D47" "46.84"
$ws.Range("E47").Value = "  +1.43%  "
Set-TextValue "D48" "1.28"
$ws.Range("E48").Value = "  -6.61%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D49" "7.59"
$ws.Range("E49").Value = "  -2.40%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D50" "2.48"
$ws.Range("E50").Value = "  -11.87%  "
Set-TextValue "D51" "0.991"
$ws.Range("E51").Value = "  -3.71%  "
